$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 187.8
$ws.Range("I42").Value = 96.333336
$ws.Range("J42").Value = 325
$ws.Range("K42").Value = 289.000008
$ws.Range("L42").Value = 975
$ws.Range("M42").Value = -59.00000799999998
$ws.Range("N42").Value = -1435

$ws.Range("H106").Value = 5005
$ws.Range("I106").Value = 5005
$ws.Range("K106").Value = 5005
$ws.Range("M106").Value = -4374

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 11584.833
$ws.Range("I45").Value = 7747.5
$ws.Range("J45").Value = 13503.5
$ws.Range("K45").Value = 7747.5
$ws.Range("L45").Value = 13503.5
$ws.Range("M45").Value = -7370.5
$ws.Range("N45").Value = -14257.5

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").ClearContents()
$ws.Range("N51").Value = 0

$ws.Range("H61").Value = 6999.2856
$ws.Range("I61").Value = 2599
$ws.Range("K61").Value = 2599
$ws.Range("M61").Value = -2387

$ws.Range("H74").Value = 9328.5
$ws.Range("I74").Value = 5323.6665
$ws.Range("K74").Value = 5323.6665
$ws.Range("M74").Value = -4449.6665

$ws.Range("H77").Value = 9328.5
$ws.Range("I77").Value = 5323.6665
$ws.Range("K77").Value = 26618.3325
$ws.Range("M77").Value = -22250.3325

$ws.Range("H136").Value = 6999.2856
$ws.Range("I136").Value = 2599
$ws.Range("K136").Value = 7797
$ws.Range("M136").Value = -5247

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1333
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1333
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = 1333
$ws.Range("N5").Value = -1559

$ws.Range("H60").Value = 99980
$ws.Range("J60").Value = 99980
$ws.Range("L60").Value = 99980
$ws.Range("N60").Value = -101178

$ws.Range("H88").Value = 18000
$ws.Range("J88").Value = 18000
$ws.Range("L88").Value = 18000
$ws.Range("N88").Value = -18812

$ws.Range("H91").Value = 18000
$ws.Range("J91").Value = 18000
$ws.Range("L91").Value = 18000
$ws.Range("N91").Value = -20808

$ws.Range("H100").Value = 18666.666
$ws.Range("J100").Value = 18666.666
$ws.Range("L100").Value = 18666.666
$ws.Range("N100").Value = -20830.666

$ws.Range("H103").Value = 22339.334
$ws.Range("J103").Value = 22339.334
$ws.Range("L103").Value = 22339.334
$ws.Range("N103").Value = -24683.334

$ws.Range("H107").Value = 1531.6
$ws.Range("I107").Value = 1531.6
$ws.Range("K107").Value = 1531.6
$ws.Range("M107").Value = 388.4000000000001

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("N120").Value = 0

$ws.Range("H134").Value = 5867.5
$ws.Range("I134").Value = 3638
$ws.Range("K134").Value = 10914
$ws.Range("M134").Value = -8379

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1255.5
$ws.Range("I16").Value = 511
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 511
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -224
$ws.Range("N16").Value = -2574

$ws.Range("H25").Value = 4000
$ws.Range("I25").Value = 4000
$ws.Range("K25").Value = 4000
$ws.Range("M25").Value = -3826

$ws.Range("H105").Value = 754.5
$ws.Range("I105").Value = 754.5
$ws.Range("K105").Value = 754.5
$ws.Range("M105").Value = 992.5

$ws.Range("H107").Value = 1205.8334
$ws.Range("I107").Value = 1011.6667
$ws.Range("K107").Value = 1011.6667
$ws.Range("M107").Value = 908.3333

$ws.Range("H113").Value = 1255.5
$ws.Range("I113").Value = 511
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 511
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1659
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 5250
$ws.Range("I11").Value = 3000
$ws.Range("K11").Value = 9000
$ws.Range("M11").Value = -8860

$ws.Range("H41").Value = 200
$ws.Range("I41").Value = 200
$ws.Range("K41").Value = 600
$ws.Range("M41").Value = -262

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 303899.25
$ws.Range("I101").Value = 1100597
$ws.Range("J101").Value = 38333.332
$ws.Range("K101").Value = 1100597
$ws.Range("L101").Value = 38333.332
$ws.Range("M101").Value = -1097352
$ws.Range("N101").Value = -44823.332

$ws.Range("H102").Value = 3130.818
$ws.Range("I102").Value = 2073.889
$ws.Range("K102").Value = 2073.889
$ws.Range("M102").Value = -451.8890000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1750
$ws.Range("I10").Value = 1750
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1750
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -1610

$ws.Range("H12").Value = 15000
$ws.Range("J12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("N12").Value = -10340

$ws.Range("H13").Value = 19999
$ws.Range("I13").Value = 19999
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 19999
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -19859

$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -2705
$ws.Range("N22").Value = -2590

$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -2893
$ws.Range("N27").Value = -2214

$ws.Range("H82").Value = 710
$ws.Range("I82").Value = 676.6
$ws.Range("K82").Value = 676.6
$ws.Range("M82").Value = -315.6

$ws.Range("H85").Value = 710
$ws.Range("I85").Value = 676.6
$ws.Range("K85").Value = 676.6
$ws.Range("M85").Value = 571.4

$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492

$ws.Range("H136").Value = 13667.667
$ws.Range("I136").Value = 9501.5
$ws.Range("K136").Value = 28504.5
$ws.Range("M136").Value = -25954.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("N27").Value = 0

$ws.Range("H68").Value = 39999
$ws.Range("J68").Value = 39999
$ws.Range("L68").Value = 39999
$ws.Range("N68").Value = -41621

$ws.Range("H71").Value = 39999
$ws.Range("J71").Value = 39999
$ws.Range("L71").Value = 119997
$ws.Range("N71").Value = -128109

$ws.Range("H101").Value = 26250
$ws.Range("J101").Value = 26250
$ws.Range("L101").Value = 26250
$ws.Range("N101").Value = -32740

$ws.Range("H104").Value = 29184.5
$ws.Range("J104").Value = 29184.5
$ws.Range("L104").Value = 29184.5
$ws.Range("N104").Value = -36172.5

$ws.Range("H136").Value = 16000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 16000
$ws.Range("K136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("M136").Value = 48000
$ws.Range("N136").Value = -53100
